# actualizacion 6 de Agosto
# Update the ICX (row 17) and NPS (row 18) indicator values on the
# "Resumen" sheet for REVAL / VALE+ / Total (columns B, C, D).
#   ICX : 4.83  -> 4.79   (entered as text, matches the original text-with-
#                           quote-prefix storage of the source cells)
#   NPS : 83.49% -> 82.33% (text value, percentage number format retained)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: ICX ---------------------------------------------------------
$ws.Range("B17").Value = "'4.79"
$ws.Range("C17").Value = "'4.79"
$ws.Range("D17").Value = "'4.79"

# --- Row 18: NPS -----------------------------------------------------------
$ws.Range("B18").Value = "'82.33%"
$ws.Range("C18").Value = "'82.33%"
$ws.Range("D18").Value = "'82.33%"
$ws.Range("B18:D18").NumberFormat = "0.00%"

# Match the saved cursor position left behind in the workbook.
$ws.Range("D18").Select() | Out-Null
